$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51; existing rows 51-193 shift down to 52-194.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new weekly data point.
$ws.Range("A51").Value = 5
$ws.Range("B51").Value = "Macroferia Regional de Talca"
$ws.Range("C51").Value = "Maule"
$ws.Range("D51").Value = 44497
$ws.Range("E51").Value = 7
$ws.Range("F51").Value = 100112009
$ws.Range("G51").Value = "Acelga"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 500
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 2000
$ws.Range("M51").Value = 2000
$ws.Range("N51").Value = "$/docena de atados (4 kilos)"
$ws.Range("O51").Value = "Región del Maule"
$ws.Range("P51").Value = 500
$ws.Range("Q51").Value = 4
$ws.Range("R51").Value = "Hortaliza"
